$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.802.38'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '2.530.02'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.575'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.51%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.534'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.88'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0811'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.62'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.00%  '
$ws.Range("E13").Value = '  -2.60%  '
$ws.Range("D14").Value = '2.911.63'
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").Value = '2.514.57'
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.14'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.849'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("D18").Value = '42.847.01'
$ws.Range("E18").Value = '  +0.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.18%  '
$ws.Range("D21").Value = '0.0₃0962'
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.66'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.95'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.45'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.60%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  +2.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '41.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.42'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.97%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '159.54'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.63%  '
$ws.Range("E33").Value = '  +3.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.71'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.32'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.89'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0790'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.113'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +16.93%  '
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '21.87'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -9.99%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.84'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0305'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.40%  '
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.47%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.022.19'
$ws.Range("E46").Value = '  -2.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '84.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '106.02'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.35'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.36%  '
$ws.Range("D51").Value = '2.768.67'
$ws.Range("E51").Value = '  +0.42%  '
